# Auto-generated cell update script
# Applies per-row market price refresh values per sheet, as scraped from the
# "Leviathan_Profits" scheduled pricing export diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H101").Value = 1000861.4
$ws.Range("I101").Value = 1667671.9
$ws.Range("J101").Value = 200688.8
$ws.Range("K101").Value = 5003015.699999999
$ws.Range("L101").Value = 602066.3999999999
$ws.Range("M101").Value = -5001393.699999999
$ws.Range("N101").Value = -605310.3999999999

$ws.Range("H107").Value = 25813.467
$ws.Range("I107").Value = 524.4783
$ws.Range("J107").Value = 108905.86
$ws.Range("K107").Value = 524.4783
$ws.Range("L107").Value = 108905.86
$ws.Range("M107").Value = 1395.5217
$ws.Range("N107").Value = -112745.86

$ws.Range("H125").Value = 13508.936
$ws.Range("I125").Value = 25063.533
$ws.Range("J125").Value = 2676.5
$ws.Range("K125").Value = 225571.797
$ws.Range("L125").Value = 24088.5
$ws.Range("M125").Value = -223111.797
$ws.Range("N125").Value = -29008.5

$ws.Range("H129").Value = 1382.2
$ws.Range("I129").Value = 552.6429000000001
$ws.Range("J129").Value = 3317.8333
$ws.Range("K129").Value = 1657.9287
$ws.Range("L129").Value = 9953.499899999999
$ws.Range("M129").Value = 3342.0713
$ws.Range("N129").Value = -19953.4999

$ws.Range("H132").Value = 1205.3429
$ws.Range("I132").Value = 1205.3429
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 3616.0287
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -1086.0287

$ws.Range("H137").Value = 2197.4644
$ws.Range("I137").Value = 2002.579
$ws.Range("K137").Value = 6007.737
$ws.Range("M137").Value = -3457.737

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20692
$ws.Range("I32").Value = 4375.873
$ws.Range("J32").Value = 84936.75
$ws.Range("K32").Value = 4375.873
$ws.Range("L32").Value = 84936.75
$ws.Range("M32").Value = -4088.873
$ws.Range("N32").Value = -85510.75

$ws.Range("H61").Value = 1681.4
$ws.Range("I61").Value = 1658.6428
$ws.Range("J61").Value = 2000
$ws.Range("K61").Value = 1658.6428
$ws.Range("L61").Value = 2000
$ws.Range("M61").Value = -1446.6428
$ws.Range("N61").Value = -2424

$ws.Range("H74").Value = 1649.2307
$ws.Range("I74").Value = 1370.0416
$ws.Range("J74").Value = 4999.5
$ws.Range("K74").Value = 1370.0416
$ws.Range("L74").Value = 4999.5
$ws.Range("M74").Value = -496.0416
$ws.Range("N74").Value = -6747.5

$ws.Range("H77").Value = 1649.2307
$ws.Range("I77").Value = 1370.0416
$ws.Range("J77").Value = 4999.5
$ws.Range("K77").Value = 6850.208000000001
$ws.Range("L77").Value = 24997.5
$ws.Range("M77").Value = -2482.208000000001
$ws.Range("N77").Value = -33733.5

$ws.Range("H136").Value = 1681.4
$ws.Range("I136").Value = 1658.6428
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 4975.928400000001
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -2425.928400000001
$ws.Range("N136").Value = -11100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 14178.066
$ws.Range("I20").Value = 15434.818
$ws.Range("J20").Value = 10722
$ws.Range("K20").Value = 15434.818
$ws.Range("L20").Value = 10722
$ws.Range("M20").Value = -15187.818
$ws.Range("N20").Value = -11216

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 23637.234
$ws.Range("I31").Value = 26482.596
$ws.Range("J31").Value = 18205.182
$ws.Range("K31").Value = 26482.596
$ws.Range("L31").Value = 18205.182
$ws.Range("M31").Value = -26187.596
$ws.Range("N31").Value = -18795.182

$ws.Range("H34").Value = 23637.234
$ws.Range("I34").Value = 26482.596
$ws.Range("J34").Value = 18205.182
$ws.Range("K34").Value = 26482.596
$ws.Range("L34").Value = 18205.182
$ws.Range("M34").Value = -26280.596
$ws.Range("N34").Value = -18609.182

$ws.Range("H134").Value = 2216.5
$ws.Range("I134").Value = 2068.9539
$ws.Range("J134").Value = 3282.111
$ws.Range("K134").Value = 6206.861699999999
$ws.Range("L134").Value = 9846.332999999999
$ws.Range("M134").Value = -3671.861699999999
$ws.Range("N134").Value = -14916.333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 1012
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()

$ws.Range("H72").Value = 1012
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H42").Value = 142645
$ws.Range("J42").Value = 142645
$ws.Range("L42").Value = 142645
$ws.Range("N42").Value = -143615

$ws.Range("H102").Value = 3665.7273
$ws.Range("I102").Value = 3632.3
$ws.Range("J102").Value = 4000
$ws.Range("K102").Value = 3632.3
$ws.Range("L102").Value = 4000
$ws.Range("M102").Value = -2010.3
$ws.Range("N102").Value = -7244

$ws.Range("H115").Value = 142645
$ws.Range("J115").Value = 142645
$ws.Range("L115").Value = 142645
$ws.Range("N115").Value = -144995

$ws.Range("H126").Value = 2148.2778
$ws.Range("I126").Value = 1780.4615
$ws.Range("J126").Value = 3104.6
$ws.Range("K126").Value = 5341.3845
$ws.Range("L126").Value = 9313.799999999999
$ws.Range("M126").Value = -2871.3845
$ws.Range("N126").Value = -14253.8

$ws.Range("H132").Value = 2930.8057
$ws.Range("I132").Value = 2263.4285
$ws.Range("K132").Value = 6790.2855
$ws.Range("M132").Value = -4260.2855

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6464.278
$ws.Range("I7").Value = 6767.2085
$ws.Range("J7").Value = 5858.4165
$ws.Range("K7").Value = 6767.2085
$ws.Range("L7").Value = 5858.4165
$ws.Range("M7").Value = -6655.2085
$ws.Range("N7").Value = -6082.4165

$ws.Range("H40").Value = 6224.7666
$ws.Range("I40").Value = 4594.5
$ws.Range("J40").Value = 8087.9287
$ws.Range("K40").Value = 4594.5
$ws.Range("L40").Value = 8087.9287
$ws.Range("M40").Value = -4458.5
$ws.Range("N40").Value = -8359.9287

$ws.Range("H93").Value = 15919.608
$ws.Range("I93").Value = 1533.3572
$ws.Range("K93").Value = 1533.3572
$ws.Range("M93").Value = -285.3571999999999

$ws.Range("H126").Value = 6464.278
$ws.Range("I126").Value = 6767.2085
$ws.Range("J126").Value = 5858.4165
$ws.Range("K126").Value = 20301.6255
$ws.Range("L126").Value = 17575.2495
$ws.Range("M126").Value = -17831.6255
$ws.Range("N126").Value = -22515.2495

$ws.Range("H136").Value = 2273.1592
$ws.Range("I136").Value = 1718.8206
$ws.Range("K136").Value = 5156.4618
$ws.Range("M136").Value = -2606.4618

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2965.3333
$ws.Range("I126").Value = 2918.6
$ws.Range("K126").Value = 8755.799999999999
$ws.Range("M126").Value = -6285.799999999999

$ws.Range("H136").Value = 792.4737
$ws.Range("I136").Value = 709.2353000000001
$ws.Range("K136").Value = 2127.7059
$ws.Range("M136").Value = 422.2941000000001
